$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '64.506.10'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -2.00%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.624.93'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -1.89%  '
$ws.Range('E4').Value = '  +0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '578.65'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -3.62%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '156.18'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.97%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.645'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +5.44%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -5.39%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '5.79'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -1.24%  '
$ws.Range('E11').Value = '  -2.70%  '
$ws.Range('E12').Value = '  +0.08%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '28.54'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('E14').Value = '  -7.12%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '3.098.38'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -1.83%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '64.372.25'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -1.96%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '2.632.40'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -1.29%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '12.26'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -4.05%  '
$ws.Range('E19').Value = '  -2.74%  '
$ws.Range('E20').Value = '  -2.12%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '345.93'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.98%  '
$ws.Range('E22').Value = '  -0.16%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '67.93'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -2.36%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.0000113'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -3.72%  '
$ws.Range('E25').Value = '  +3.71%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '9.37'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -3.83%  '
$ws.Range('E27').Value = '  -2.82%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '555.30'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +4.05%  '
$ws.Range('E29').Value = '  -2.33%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '7.98'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.77%  '
$ws.Range('E31').Value = '  +0.01%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '2.08'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -2.74%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.73'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -2.53%  '
$ws.Range('E34').Value = '  -1.20%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '5.35'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -2.50%  '
$ws.Range('E36').Value = '  -2.74%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '20.02'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -3.17%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('E39').Value = '  -0.29%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '151.45'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -4.22%  '
$ws.Range('E41').Value = '  -0.04%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.46'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +3.30%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '158.21'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -3.21%  '
$ws.Range('E44').Value = '  -2.96%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0601'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -2.52%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '22.91'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +0.46%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.634'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.15%  '
$ws.Range('E49').Value = '  -3.75%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '19.17'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -4.73%  '
$ws.Range('E51').Value = '  -6.34%  '
